$d = $word.ActiveDocument

# 1) "Don't delay! Book your spot today!" -> Russian (unique occurrence)
$d.Content.Find.Execute(
    "Don’t delay! Book your spot today!", $true, $false, $false, $false, $false,
    $true, 1, $false, "Не откладывайте! Забронируйте свое место сегодня!", 2
) | Out-Null

# 2) "We look forward to seeing you at [EVENT NAME]! " appears twice with identical
#    visible text, but only the SECOND occurrence (a single plain run, matching the
#    diff's single <w:r> replacement) should be translated. The first occurrence is
#    split across three runs (with "[EVENT NAME]" highlighted) and must stay in English.
#    Find the first occurrence's paragraph, then restrict the search range to start
#    right after it so only the second occurrence gets replaced.
$firstOccurrence = $d.Content
$firstOccurrence.Find.Execute("We look forward to seeing you at [EVENT NAME]! ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterFirst = $d.Range($firstOccurrence.End, $d.Content.End)
$afterFirst.Find.Execute(
    "We look forward to seeing you at [EVENT NAME]! ", $true, $false, $false, $false, $false,
    $true, 1, $false, "С нетерпением ждем встречи на [EVENT NAME]! ", 2
) | Out-Null

# 3) "If you have any questions, please contact your country manager:" -> Russian (unique occurrence)
$d.Content.Find.Execute(
    "If you have any questions, please contact your country manager:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Если у вас возникли вопросы, обратитесь к вашему региональному менеджеру:", 2
) | Out-Null

# 4) "If you have any questions, please contact us via:" -> Russian (unique occurrence)
$d.Content.Find.Execute(
    "If you have any questions, please contact us via:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Свяжитесь с нами, если у вас есть вопросы:", 2
) | Out-Null
